$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2 = 92.01441192626953
    3 = 92.06435394287109
    4 = 94.12319183349609
    5 = 94.20018005371094
    6 = 91.93358612060547
    7 = 92.01475524902344
    8 = 90.50289154052734
    9 = 90.61515045166016
    10 = 89.81863403320312
    11 = 89.92295074462891
    12 = 89.08560943603516
    13 = 89.18537139892578
    14 = 96.22895812988281
    15 = 96.32610321044922
    16 = 114.5813369750977
    17 = 114.6759948730469
    18 = 128.3124542236328
    19 = 128.4047241210938
    20 = 141.5168304443359
    21 = 141.6029205322266
    22 = 141.3161926269531
    23 = 141.3960876464844
    24 = 133.8856048583984
    25 = 133.9608917236328
    26 = 136.4703979492188
    27 = 136.5442047119141
    28 = 131.8370971679688
    29 = 131.9131927490234
    30 = 136.4724426269531
    31 = 136.5547637939453
    32 = 150.7989654541016
    33 = 150.8905792236328
    34 = 178.1812286376953
    35 = 178.2843017578125
    36 = 202.3189544677734
    37 = 202.4351654052734
    38 = 174.9986267089844
    39 = 175.1290740966797
    40 = 138.7809600830078
    41 = 138.9259185791016
    42 = 120.2549209594727
    43 = 120.4136581420898
    44 = 106.5871047973633
    45 = 106.7581939697266
    46 = 94.4610595703125
    47 = 94.64311218261719
    48 = 87.37326049804688
    49 = 87.56507873535156
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 2).Value = $newValues[$row]
}
